# Applies the "Fixed errors for demo. Added eCommerceSelected intent" edit
# to the UserSays / entities worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- F9: new "Puedes comprarlo aqui..." help text (was empty) ---------------
$ws.Range("F9").Value = "Puedes comprarlo aquí`nEstas son las tiendas donde puedes comprarlo`nAquí tienes las tiendas online donde lo ofrecen"
$ws.Range("F9").WrapText = $true
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("F9").VerticalAlignment = -4108

# --- F23: reworded "quick charge" explanation --------------------------------
$ws.Range("F23").Value = "La carga rápdia por ahora solo suele venir en móviles de gama alta o premium. Aquí los tienes."

# --- F26: same reworded text now also shown here (was empty) ----------------
$ws.Range("F26").Value = "La carga rápdia por ahora solo suele venir en móviles de gama alta o premium. Aquí los tienes."

# --- F27: reworded "ask or help" prompt --------------------------------------
$ws.Range("F27").Value = "Puedes preguntarme lo que quieras sobre móviles, o te puedo guiar para elegir uno.`n ¿Quieres preguntarme algo? ¿O prefieres que te ayude a elegir características?`n"

# --- Row 28: brand-new "eCommerceSelected" intent row ------------------------
$ws.Range("A28").Value = "sp.selected.ecommerce.selected"
$ws.Range("B28").Value = "Amazon`nMediaMarkt`nAliExpress`nGearBest"
$ws.Range("C28").Value = "eCommerceName"
$ws.Range("D28").Value = "Amazon`nMediaMarkt`nAliExpress`nGearBest"
$ws.Range("F28").Value = "Aquí tienes los detalles del sitio de compra elegido"

$ws.Range("B28:D28").WrapText = $true
$ws.Range("A28:F28").HorizontalAlignment = -4108
$ws.Range("A28:F28").VerticalAlignment = -4108

$ws.Rows.Item(28).RowHeight = 70.5

# --- Selection / view bookkeeping (cosmetic, matches the saved view state) --
$ws.Range("D32").Select()
